# Refresh the "cryptos" price/volume snapshot (GitHub Actions scrape update).
# Most rows only get new Price (D) / Volume(1h) (E) text.
# Rows 36-37 additionally swap Monero/NEARProtocol into new rank order, so
# their Coin (B) and Link (C) columns are rewritten too.
#
# Price values that happen to parse as a plain number (e.g. "552.32",
# "0.0000163") are forced to Text format first, matching the source data
# which stores every Price/Volume cell as a string (inline string), not a
# numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.784.94"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "3.276.62"
$ws.Range("E3").Value = "  -4.03%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.32"
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.39"
$ws.Range("E6").Value = "  -8.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.278.38"
$ws.Range("E8").Value = "  -3.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.78"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  -5.19%  "
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("D13").Value = "3.833.42"
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.58"
$ws.Range("E15").Value = "  -7.35%  "
$ws.Range("D16").Value = "3.274.11"
$ws.Range("E16").Value = "  -4.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000163"
$ws.Range("E17").Value = "  -5.35%  "
$ws.Range("D18").Value = "59.896.05"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("E19").Value = "  -6.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.68"
$ws.Range("E20").Value = "  -5.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.51"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.18"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.34"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("E25").Value = "  -7.31%  "
$ws.Range("D26").Value = "3.411.50"
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000101"
$ws.Range("E27").Value = "  -9.97%  "
$ws.Range("E28").Value = "  -5.46%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("E30").Value = "  -8.53%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.44"
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.46"
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  -7.76%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.05"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.03"
$ws.Range("E37").Value = "  -8.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  -6.04%  "
$ws.Range("E39").Value = "  -5.22%  "
$ws.Range("D40").Value = "3.302.46"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.96"
$ws.Range("E41").Value = "  -16.16%  "
$ws.Range("E42").Value = "  -7.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.67"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.744"
$ws.Range("E44").Value = "  -4.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.10"
$ws.Range("E45").Value = "  -7.28%  "
$ws.Range("E46").Value = "  -7.19%  "
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "2.323.83"
$ws.Range("E49").Value = "  -8.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.32"
$ws.Range("E50").Value = "  -7.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.11"
$ws.Range("E51").Value = "  -6.66%  "
